$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9198808073997498
$ws.Range("B1").Value = 1.76668381690979
$ws.Range("C1").Value = 8.63957405090332
$ws.Range("D1").Value = 2.020910978317261
$ws.Range("E1").Value = 1.185853600502014
